$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 9157.286
$ws.Range("I18").Value = 10150.25
$ws.Range("K18").Value = 10150.25
$ws.Range("M18").Value = -9866.25

$ws.Range("H51").Value = 2875
$ws.Range("J51").Value = 2875
$ws.Range("L51").Value = 2875
$ws.Range("N51").Value = -3843

$ws.Range("H87").Value = 84657.164
$ws.Range("J87").Value = 94588.60000000001
$ws.Range("L87").Value = 94588.60000000001
$ws.Range("N87").Value = -97084.60000000001

$ws.Range("H90").Value = 84657.164
$ws.Range("J90").Value = 94588.60000000001
$ws.Range("L90").Value = 283765.8
$ws.Range("N90").Value = -296245.8

$ws.Range("H112").Value = 1667.3125
$ws.Range("I112").Value = 821.75
$ws.Range("J112").Value = 1949.1666
$ws.Range("K112").Value = 2465.25
$ws.Range("L112").Value = 5847.4998
$ws.Range("M112").Value = -1357.25
$ws.Range("N112").Value = -8063.4998

$ws.Range("H121").Value = 1275
$ws.Range("J121").Value = 1275
$ws.Range("L121").Value = 3825
$ws.Range("N121").Value = -7319

$ws.Range("H125").Value = 41668580
$ws.Range("J125").Value = 55557428
$ws.Range("L125").Value = 500016852
$ws.Range("N125").Value = -500021772

$ws.Range("H132").Value = 29130.74
$ws.Range("I132").Value = 31901.658
$ws.Range("K132").Value = 95704.974
$ws.Range("M132").Value = -93174.974

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1026.3469
$ws.Range("I2").Value = 902.0217
$ws.Range("K2").Value = 902.0217
$ws.Range("M2").Value = -789.0217

$ws.Range("H36").Value = 800
$ws.Range("I36").Value = 800
$ws.Range("K36").Value = 800
$ws.Range("M36").Value = -454

$ws.Range("H116").Value = 1026.3469
$ws.Range("I116").Value = 902.0217
$ws.Range("K116").Value = 902.0217
$ws.Range("M116").Value = 1391.9783

$ws.Range("H122").Value = 2418.1
$ws.Range("I122").Value = 1311.6428
$ws.Range("K122").Value = 3934.9284
$ws.Range("M122").Value = -1484.9284

$ws.Range("H132").Value = 478567.22
$ws.Range("I132").Value = 666981.3
$ws.Range("K132").Value = 2000943.9
$ws.Range("M132").Value = -1998413.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1026.3469
$ws.Range("I3").Value = 902.0217
$ws.Range("K3").Value = 902.0217
$ws.Range("M3").Value = -788.0217

$ws.Range("H105").Value = 4419.346
$ws.Range("I105").Value = 3991.3
$ws.Range("K105").Value = 3991.3
$ws.Range("M105").Value = -2244.3

$ws.Range("H134").Value = 571221.4
$ws.Range("I134").Value = 704067.75
$ws.Range("J134").Value = 6624.125
$ws.Range("K134").Value = 2112203.25
$ws.Range("L134").Value = 19872.375
$ws.Range("M134").Value = -2109668.25
$ws.Range("N134").Value = -24942.375

$ws.Range("H139").Value = 54999.5
$ws.Range("J139").Value = 54999.5
$ws.Range("L139").Value = 54999.5
$ws.Range("N139").Value = -65279.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3913.8333
$ws.Range("I99").Value = 3583.375
$ws.Range("K99").Value = 3583.375
$ws.Range("M99").Value = -2085.375

$ws.Range("H126").Value = 3913.8333
$ws.Range("I126").Value = 3583.375
$ws.Range("K126").Value = 10750.125
$ws.Range("M126").Value = -8280.125

$ws.Range("H132").Value = 11381704
$ws.Range("I132").Value = 22587.941
$ws.Range("J132").Value = 50002700
$ws.Range("K132").Value = 67763.823
$ws.Range("L132").Value = 150008100
$ws.Range("M132").Value = -65233.823
$ws.Range("N132").Value = -150013160

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 888.1111
$ws.Range("I86").Value = 913.5714
$ws.Range("J86").Value = 799
$ws.Range("K86").Value = 2740.7142
$ws.Range("L86").Value = 2397
$ws.Range("M86").Value = -1554.7142
$ws.Range("N86").Value = -4769

$ws.Range("H89").Value = 888.1111
$ws.Range("I89").Value = 913.5714
$ws.Range("J89").Value = 799
$ws.Range("K89").Value = 8222.142600000001
$ws.Range("L89").Value = 7191
$ws.Range("M89").Value = -2294.142600000001
$ws.Range("N89").Value = -19047

$ws.Range("H92").Value = 663.5
$ws.Range("J92").Value = 203.75
$ws.Range("L92").Value = 611.25
$ws.Range("N92").Value = -3107.25

$ws.Range("H93").Value = 6307.3887
$ws.Range("I93").Value = 4000
$ws.Range("J93").Value = 6443.1177
$ws.Range("K93").Value = 12000
$ws.Range("L93").Value = 19329.3531
$ws.Range("M93").Value = -10128
$ws.Range("N93").Value = -23073.3531

$ws.Range("H139").Value = 691.75
$ws.Range("I139").Value = 654.8570999999999
$ws.Range("J139").Value = 950
$ws.Range("K139").Value = 1964.5713
$ws.Range("L139").Value = 2850
$ws.Range("M139").Value = 3175.4287
$ws.Range("N139").Value = -13130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 610352.9
$ws.Range("I21").Value = 5000499.5
$ws.Range("K21").Value = 5000499.5
$ws.Range("M21").Value = -5000326.5

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null

$ws.Range("H30").Value = 610352.9
$ws.Range("I30").Value = 5000499.5
$ws.Range("K30").Value = 5000499.5
$ws.Range("M30").Value = -5000394.5

$ws.Range("H70").Value = 4422.7144
$ws.Range("I70").Value = 4327.3335
$ws.Range("J70").Value = 4995
$ws.Range("K70").Value = 4327.3335
$ws.Range("L70").Value = 4995
$ws.Range("M70").Value = -4057.3335
$ws.Range("N70").Value = -5535

$ws.Range("H73").Value = 4422.7144
$ws.Range("I73").Value = 4327.3335
$ws.Range("J73").Value = 4995
$ws.Range("K73").Value = 4327.3335
$ws.Range("L73").Value = 4995
$ws.Range("M73").Value = -3391.3335
$ws.Range("N73").Value = -6867

$ws.Range("H102").Value = 2295.125
$ws.Range("I102").Value = 1591.6471
$ws.Range("K102").Value = 1591.6471
$ws.Range("M102").Value = 30.35290000000009

$ws.Range("H117").Value = 59326.75
$ws.Range("J117").Value = 59326.75
$ws.Range("L117").Value = 59326.75
$ws.Range("N117").Value = -66210.75

$ws.Range("H122").Value = 4091.4285
$ws.Range("I122").Value = 2325.8462
$ws.Range("J122").Value = 6960.5
$ws.Range("K122").Value = 6977.5386
$ws.Range("L122").Value = 20881.5
$ws.Range("M122").Value = -4527.5386
$ws.Range("N122").Value = -25781.5

$ws.Range("H126").Value = 621313.5600000001
$ws.Range("I126").Value = 1391081.5
$ws.Range("K126").Value = 4173244.5
$ws.Range("M126").Value = -4170774.5

$ws.Range("H132").Value = 1253.7188
$ws.Range("I132").Value = 1110.3405
$ws.Range("K132").Value = 3331.0215
$ws.Range("M132").Value = -801.0214999999998

$ws.Range("H134").Value = 45860.59
$ws.Range("J134").Value = 45860.59
$ws.Range("L134").Value = 137581.77
$ws.Range("N134").Value = -142651.77

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5991.4707
$ws.Range("I7").Value = 5525.357
$ws.Range("K7").Value = 5525.357
$ws.Range("M7").Value = -5413.357

$ws.Range("H22").Value = 79818.30499999999
$ws.Range("J22").Value = 3285.2727
$ws.Range("L22").Value = 3285.2727
$ws.Range("N22").Value = -3875.2727

$ws.Range("H27").Value = 79818.30499999999
$ws.Range("J27").Value = 3285.2727
$ws.Range("L27").Value = 3285.2727
$ws.Range("N27").Value = -3499.2727

$ws.Range("H38").Value = 15497.5

$ws.Range("H40").Value = 19524.5
$ws.Range("I40").Value = 19524.5
$ws.Range("K40").Value = 19524.5
$ws.Range("M40").Value = -19388.5

$ws.Range("H122").Value = 3794.0789
$ws.Range("J122").Value = 4822.5
$ws.Range("L122").Value = 14467.5
$ws.Range("N122").Value = -19367.5

$ws.Range("H126").Value = 5991.4707
$ws.Range("I126").Value = 5525.357
$ws.Range("K126").Value = 16576.071
$ws.Range("M126").Value = -14106.071

$ws.Range("H132").Value = 2879.6924
$ws.Range("I132").Value = 2228.5957
$ws.Range("K132").Value = 6685.7871
$ws.Range("M132").Value = -4155.7871

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 18000
$ws.Range("J26").Value = 17500
$ws.Range("L26").Value = 17500
$ws.Range("N26").Value = -18086

$ws.Range("H46").Value = 61917.09
$ws.Range("J46").Value = 65108.8
$ws.Range("L46").Value = 65108.8
$ws.Range("N46").Value = -65570.8

$ws.Range("H122").Value = 3098.5
$ws.Range("I122").Value = 2231.0833
$ws.Range("K122").Value = 6693.249899999999
$ws.Range("M122").Value = -4243.249899999999

$ws.Range("H132").Value = 3369.6667
$ws.Range("I132").Value = 3123.0303
$ws.Range("J132").Value = 3912.2666
$ws.Range("K132").Value = 9369.090899999999
$ws.Range("L132").Value = 11736.7998
$ws.Range("M132").Value = -6839.090899999999
$ws.Range("N132").Value = -16796.7998

$ws.Range("H134").Value = 61917.09
$ws.Range("J134").Value = 65108.8
$ws.Range("L134").Value = 195326.4
$ws.Range("N134").Value = -200396.4

$ws.Range("H136").Value = 9600731
$ws.Range("I136").Value = 10799620
$ws.Range("K136").Value = 32398860
$ws.Range("M136").Value = -32396310
